# added sorting facilities to transaction MVC
# Fill in the newly-logged Manhours / Feature / Comment entries for rows 37-41
# on the Reporting sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37 - "Add date to transactions table..." entry already has its Feature
# comment (E37); just add the Manhours and the follow-up Comment.
$ws.Range("B37").Value = 1
$ws.Range("F37").Value = "straightforward enough"

# Row 38 - new Feature / Comment pair.
$ws.Range("B38").Value = 1
$ws.Range("E38").Value = "Create struct to hold budget in tag"
$ws.Range("F38").Value = "should have continued with this"

# Row 39 - new Feature / Comment pair.
$ws.Range("B39").Value = 1
$ws.Range("E39").Value = "Display total amount spent per tag per month"
$ws.Range("F39").Value = "ERROR not due to struct but using post rather than get"

# Row 40 - new Feature / Comment pair.
$ws.Range("B40").Value = 2
$ws.Range("E40").Value = "Rebuild a budget class"
$ws.Range("F40").Value = "feels a bit surplus and details not persisted? Would have liked the struct route"

# Row 41 - new Feature / Comment pair.
$ws.Range("B41").Value = 1
$ws.Range("E41").Value = "show totals for month"
$ws.Range("F41").Value = "straight forward enough"

# Move the view's scroll position / active selection to reflect where work
# left off.
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("D42").Select()
